$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 14567441.71
$ws.Range("P2").Value = 58.9112450723
$ws.Range("Q2").Value = 132449936.58
$ws.Range("R2").Value = 535.6321877931
$ws.Range("S2").Value = 41979494.95
$ws.Range("T2").Value = 169.7665495592
$ws.Range("U2").Value = -82594731.76000001
$ws.Range("V2").Value = -334.0159913634
$ws.Range("W2").Value = 1458305.23
$ws.Range("X2").Value = 5.8974374846
$ws.Range("Y2").Value = 40155490.3
$ws.Range("Z2").Value = 162.3902107971
$ws.Range("AA2").Value = 43064880.26
$ws.Range("AB2").Value = 174.1558858111
$ws.Range("AC2").Value = -24727777.68
$ws.Range("AD2").ClearContents()
